$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7 (shifts old rows 7-21 down to rows 8-22, and
# Excel auto-adjusts the formulas that referenced those rows).
$ws.Rows("7:7").Insert()

# New row 7 picks up the date-formatted style used by the rows above it.
$ws.Range("B6").Copy($ws.Range("B7"))

# Populate the newly inserted row with the new bootcamp entry.
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = 44323
$ws.Range("C7").Value = "76: Pseudo Elements"
$ws.Range("E7").Formula = "=27+68"

# Rows 8 and 9 (formerly blank rows 7 and 8, now shifted down) just get the
# week-number filled in.
$ws.Range("A8").Value = 3
$ws.Range("A9").Value = 3

# The week boundaries shifted because of the newly inserted row, so the
# weekly roll-up formulas need to be corrected by hand to match the
# reorganized week cutoffs.
$ws.Range("E15").Formula = "=SUM(E3:E5)"
$ws.Range("E17").Formula = "=SUM(E6:E8)"

# Restore the previously-selected cell shown in the saved file.
$ws.Range("C21").Select()

$wb.Save()
